$d = $word.ActiveDocument

$replacements = @(
    @("N = 34,218", "N = 28,001"),
    @("96.4 (89.7, 102.8)", "96.2 (89.5, 102.7)"),
    @("34,100 (100)", "27,897 (100)"),
    @("33,974 (99)", "27,773 (99)"),
    @("466.8 (323.3, 641.4)", "421.1 (295.8, 567.1)"),
    @("236.1 (117.0, 410.0)", "223.5 (109.5, 387.6)"),
    @("589.5 (375.4, 856.6)", "518.7 (336.8, 728.2)"),
    @("773 (2.3)", "689 (2.5)"),
    @("497 (1.5)", "439 (1.6)"),
    @("63.6 (56.6, 68.6)", "64.0 (56.9, 68.9)"),
    @("1,025 (3.0)", "786 (2.8)"),
    @("33,193 (97)", "27,215 (97)"),
    @("19,648 (57)", "14,895 (53)"),
    @("14,570 (43)", "13,106 (47)"),
    @("2,700 (7.9)", "2,177 (7.8)"),
    @("8,281 (24)", "6,516 (23)"),
    @("8,044 (24)", "6,557 (23)"),
    @("15,193 (44)", "12,751 (46)"),
    @("4,418 (13)", "3,573 (13)"),
    @("7,562 (22)", "6,148 (22)"),
    @("8,952 (26)", "7,304 (26)"),
    @("7,825 (23)", "6,594 (24)"),
    @("2,248 (6.6)", "1,920 (6.9)"),
    @("3,213 (9.4)", "2,462 (8.8)"),
    @("197 (0.6)", "136 (0.5)"),
    @("13,588 (40)", "10,312 (37)"),
    @("14,175 (41)", "12,046 (43)"),
    @("6,258 (18)", "5,507 (20)"),
    @("33,086 (97)", "26,982 (96)"),
    @("6,824 (21)", "5,534 (20)"),
    @("26,451 (79)", "21,716 (80)"),
    @("19,710 (58)", "15,997 (57)"),
    @("12,213 (36)", "10,098 (36)"),
    @("2,295 (6.7)", "1,906 (6.8)"),
    @("1,879 (5.5)", "1,496 (5.3)"),
    @("6,826 (20)", "5,526 (20)"),
    @("8,456 (25)", "6,945 (25)"),
    @("8,987 (26)", "7,366 (26)"),
    @("8,070 (24)", "6,668 (24)"),
    @("24,361 (71)", "19,652 (70)"),
    @("8,664 (25)", "7,335 (26)"),
    @("1,193 (3.5)", "1,014 (3.6)"),
    @("6,177 (18)", "4,778 (17)"),
    @("11,667 (34)", "9,617 (34)"),
    @("16,374 (48)", "13,606 (49)"),
    @("28,275 (83)", "23,156 (83)"),
    @("5,675 (17)", "4,626 (17)"),
    @("268 (0.8)", "219 (0.8)"),
    @("29,083 (85)", "23,840 (85)"),
    @("4,979 (15)", "4,037 (14)"),
    @("156 (0.5)", "124 (0.4)"),
    @("11,115 (32)", "9,401 (34)"),
    @("20,568 (60)", "16,617 (59)"),
    @("2,535 (7.4)", "1,983 (7.1)"),
    @("5,761 (17)", "4,821 (17)"),
    @("25,213 (74)", "20,560 (73)"),
    @("3,244 (9.5)", "2,620 (9.4)"),
    @("29.9 (16.8, 50.6)", "29.8 (16.7, 50.6)"),
    @("19.6 (6.0, 41.5)", "19.4 (5.8, 41.5)")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $old"
    }
}

$d.Save()
